$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Summary block (rows 10-12): plug in the real right/wrong/not-attempt/marks
# numbers now that the marking can cope with a float/negative entry instead
# of silently falling back to the "Absent" placeholder.
# ---------------------------------------------------------------------------

# Row 10 ("No.") - counts
$ws.Range("B10").Value = 22
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 28

# Row 11 ("Marking") - per-question marks; C11 becomes a real number (-1)
# instead of the text string "-1" that used to trip up downstream reads.
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 ("Total")
$ws.Range("B12").Value = 88
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "87/112"

# A10/A11/A12 pick up the bold "mtitleStyle" formatting (style index 4) that
# the rest of the header row (A9:E9) already uses.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Answer grid: the sheet used to show three side-by-side Student/Correct Ans
# blocks (A:B, D:E, G:H) but only ever populated the first 3-5 rows of the
# second/third blocks. Now it's a single two-column block (A:B) that spans
# every question, plus the small 3-row leftover block (D:E) for rows 16-18.
# The G:H block is dropped entirely.
# ---------------------------------------------------------------------------

# Fill in the student's answer (column A) for every attempted question,
# rows 16-40, matching column B (the correct answer) when right.
$ws.Range("A16").Value = "Option A"
$ws.Range("A18").Value = "Option B"
$ws.Range("A20").Value = "Option B"
$ws.Range("A22").Value = "Option D"
$ws.Range("A23").Value = "Option D"
$ws.Range("A26").Value = "Option D"
$ws.Range("A27").Value = "Option A"
$ws.Range("A28").Value = "Option D"
$ws.Range("A29").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("A31").Value = "Option D"
$ws.Range("A32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("A34").Value = "Option B"
$ws.Range("A35").Value = "Option D"
$ws.Range("A36").Value = "Option A"
$ws.Range("A37").Value = "Option A"
$ws.Range("A38").Value = "Option A"
$ws.Range("A39").Value = "Option D"
$ws.Range("A40").Value = "Option D"
# Rows 17, 19, 21, 24, 25 stay blank (not attempted).

# Style the newly-filled student answers: green ("correctStyle", same as the
# existing B column) where it matches the correct answer, red
# ("incorrectStyle") where it doesn't (only row 26: answered D, correct C).
$correctRows = @(16, 18, 20, 22, 23, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40)
foreach ($r in $correctRows) {
    $ws.Range("B10").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
}
$ws.Range("C10").Copy()
$ws.Range("A26").PasteSpecial(-4122)

# D16:E18 leftover block keeps its 3 rows; fill in the student-answer column
# (D) to mirror the correct answer (E) since all three were right.
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option D"
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("D18").PasteSpecial(-4122)

# Drop the rest of the D:E block (rows 19-40) and the whole G:H block
# (rows 15-40) - they're no longer part of the sheet's used range.
$ws.Range("D19:E40").Clear()
$ws.Range("G15:H40").Clear()
